# Data update from DGS's 2021/09/25, 2021/09/26 and 2021/09/27 reports.
# Append the latest report row (row 86) to the risk-matrix time series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 86

# Column A holds the date label as *text* (the whole column is authored
# that way even though it is styled with a yyyy/mm/dd number format), so
# force text entry to stop Excel from auto-converting the "yyyy/mm/dd"
# looking string into a date serial number, then restore the column's
# usual number format afterwards.
$aCell = $ws.Cells.Item($newRow, 1)
$aCell.NumberFormat = "@"
$aCell.Value = "2021/09/27"
$aCell.NumberFormat = "yyyy/mm/dd"

$ws.Cells.Item($newRow, 2).Value = 111.6
$ws.Cells.Item($newRow, 3).Value = 113.5
$ws.Cells.Item($newRow, 4).Value = 0.85
$ws.Cells.Item($newRow, 5).Value = 0.84

# Put the selection on the next empty row, as Excel does after data entry.
$ws.Range("A87").Select() | Out-Null
